$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4
$ws.Range("A4").Value = -21.756
$ws.Range("D4").Value = -7.833

# Row 5
$ws.Range("D5").Value = -8.206999999999999

# Row 7
$ws.Range("A7").Value = -20.987

# Row 8
$ws.Range("D8").Value = -7.896000000000001

# Row 16
$ws.Range("A16").Value = -20.654
$ws.Range("D16").Value = -8.463999999999999
